$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 100068.4
$ws.Range("I6").Value = 125061.75
$ws.Range("K6").Value = 375185.25
$ws.Range("M6").Value = -375073.25

$ws.Range("H20").Value = 4055.5
$ws.Range("I20").Value = 1111
$ws.Range("J20").Value = 7000
$ws.Range("K20").Value = 1111
$ws.Range("L20").Value = 7000
$ws.Range("M20").Value = -881
$ws.Range("N20").Value = -7460

$ws.Range("H32").Value = 2782
$ws.Range("I32").Value = 5635
$ws.Range("J32").Value = 1831
$ws.Range("K32").Value = 5635
$ws.Range("L32").Value = 1831
$ws.Range("M32").Value = -5309
$ws.Range("N32").Value = -2483

$ws.Range("H33").Value = 747.125
$ws.Range("J33").Value = 754.5
$ws.Range("L33").Value = 754.5
$ws.Range("N33").Value = -1212.5

$ws.Range("H35").Value = 4055.5
$ws.Range("I35").Value = 1111
$ws.Range("J35").Value = 7000
$ws.Range("K35").Value = 1111
$ws.Range("L35").Value = 7000
$ws.Range("M35").Value = -732
$ws.Range("N35").Value = -7758

$ws.Range("H138").Value = 3488.7778
$ws.Range("I138").Value = 2995.25
$ws.Range("J138").Value = 3574.6086
$ws.Range("K138").Value = 8985.75
$ws.Range("L138").Value = 10723.8258
$ws.Range("M138").Value = -3845.75
$ws.Range("N138").Value = -21003.8258

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3213.5454
$ws.Range("I2").Value = 2261
$ws.Range("K2").Value = 2261
$ws.Range("M2").Value = -2148

$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H32").Value = 7507.5366
$ws.Range("I32").Value = 7507.5366
$ws.Range("K32").Value = 7507.5366
$ws.Range("M32").Value = -7220.5366

$ws.Range("H41").Value = 30475
$ws.Range("I41").Value = 1150
$ws.Range("J41").Value = 34986.54
$ws.Range("K41").Value = 1150
$ws.Range("L41").Value = 34986.54
$ws.Range("M41").Value = -736
$ws.Range("N41").Value = -35814.54

$ws.Range("H46").Value = 7821.75
$ws.Range("J46").Value = 8346.444
$ws.Range("L46").Value = 8346.444
$ws.Range("N46").Value = -8984.444

$ws.Range("H74").Value = 4241.143
$ws.Range("I74").Value = 3124.4783
$ws.Range("K74").Value = 3124.4783
$ws.Range("M74").Value = -2250.4783

$ws.Range("H77").Value = 4241.143
$ws.Range("I77").Value = 3124.4783
$ws.Range("K77").Value = 15622.3915
$ws.Range("M77").Value = -11254.3915

$ws.Range("H110").Value = 6598.933
$ws.Range("I110").Value = 4599.8
$ws.Range("K110").Value = 4599.8
$ws.Range("M110").Value = -2554.8

$ws.Range("H116").Value = 3213.5454
$ws.Range("I116").Value = 2261
$ws.Range("K116").Value = 2261
$ws.Range("M116").Value = 33

$ws.Range("H132").Value = 1612.1428
$ws.Range("I132").Value = 1571.3158
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 4713.9474
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -2183.9474
$ws.Range("N132").Value = -11060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3213.5454
$ws.Range("I3").Value = 2261
$ws.Range("K3").Value = 2261
$ws.Range("M3").Value = -2147

$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H105").Value = 2627.35
$ws.Range("I105").Value = 2682.1333
$ws.Range("K105").Value = 2682.1333
$ws.Range("M105").Value = -935.1333

$ws.Range("H134").Value = 6248.75
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 6248.75
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 18746.25
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -23816.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3252.1724
$ws.Range("I16").Value = 2101.7368
$ws.Range("K16").Value = 2101.7368
$ws.Range("M16").Value = -1814.7368

$ws.Range("H22").Value = 912.6
$ws.Range("I22").Value = 964.2143
$ws.Range("J22").Value = 190
$ws.Range("K22").Value = 964.2143
$ws.Range("L22").Value = 190
$ws.Range("M22").Value = -614.2143
$ws.Range("N22").Value = -890

$ws.Range("H31").Value = 33339278
$ws.Range("I31").Value = 62503224
$ws.Range("K31").Value = 62503224
$ws.Range("M31").Value = -62502929

$ws.Range("H34").Value = 33339278
$ws.Range("I34").Value = 62503224
$ws.Range("K34").Value = 62503224
$ws.Range("M34").Value = -62503022

$ws.Range("H113").Value = 3252.1724
$ws.Range("I113").Value = 2101.7368
$ws.Range("K113").Value = 2101.7368
$ws.Range("M113").Value = 68.26319999999987

$ws.Range("H141").Value = 193944.84
$ws.Range("J141").Value = 225526.9
$ws.Range("L141").Value = 225526.9
$ws.Range("N141").Value = -235886.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2856.4285
$ws.Range("I114").Value = 999
$ws.Range("J114").Value = 3599.4
$ws.Range("K114").Value = 2997
$ws.Range("L114").Value = 10798.2
$ws.Range("M114").Value = 257
$ws.Range("N114").Value = -17306.2

$ws.Range("H121").Value = 6313501
$ws.Range("I121").Value = 373.375
$ws.Range("J121").Value = 12626629
$ws.Range("K121").Value = 1120.125
$ws.Range("L121").Value = 37879887
$ws.Range("M121").Value = 189.875
$ws.Range("N121").Value = -37882507

$ws.Range("H132").Value = 31252010
$ws.Range("I132").Value = 62501520
$ws.Range("K132").Value = 562513680
$ws.Range("M132").Value = -562511150

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2020.5264
$ws.Range("I132").Value = 2020.5264
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6061.5792
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3531.5792
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1554.3636
$ws.Range("I22").Value = 1549.6666
$ws.Range("J22").Value = 1560
$ws.Range("K22").Value = 1549.6666
$ws.Range("L22").Value = 1560
$ws.Range("M22").Value = -1254.6666
$ws.Range("N22").Value = -2150

$ws.Range("H27").Value = 1554.3636
$ws.Range("I27").Value = 1549.6666
$ws.Range("J27").Value = 1560
$ws.Range("K27").Value = 1549.6666
$ws.Range("L27").Value = 1560
$ws.Range("M27").Value = -1442.6666
$ws.Range("N27").Value = -1774

$ws.Range("H93").Value = 16460.143
$ws.Range("I93").Value = 2355.1
$ws.Range("K93").Value = 2355.1
$ws.Range("M93").Value = -1107.1

$ws.Range("H100").Value = 7989.8
$ws.Range("I100").Value = 7474.5
$ws.Range("J100").Value = 8333.333000000001
$ws.Range("K100").Value = 7474.5
$ws.Range("L100").Value = 8333.333000000001
$ws.Range("M100").Value = -6933.5
$ws.Range("N100").Value = -9415.333000000001

$ws.Range("I136").Value = 55564908
$ws.Range("K136").Value = 166694724
$ws.Range("M136").Value = -166692174

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1722.9231
$ws.Range("I2").Value = 924.75
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 924.75
$ws.Range("L2").Value = 3000
$ws.Range("M2").Value = -812.75
$ws.Range("N2").Value = -3224

$ws.Range("H122").Value = 3980.4119
$ws.Range("I122").Value = 2885.2
$ws.Range("K122").Value = 8655.599999999999
$ws.Range("M122").Value = -6205.599999999999

$ws.Range("H132").Value = 4679.385
$ws.Range("I132").Value = 3713.2646
$ws.Range("K132").Value = 11139.7938
$ws.Range("M132").Value = -8609.793799999999

$ws.Range("H136").Value = 1835.2565
$ws.Range("I136").Value = 1173.6
$ws.Range("K136").Value = 3520.8
$ws.Range("M136").Value = -970.7999999999997
